# Updated symbol list on Tue Jan 17 03:42:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.41"
$ws.Range("E2").Value = "'-1.72%"
$ws.Range("D3").Value = "'31.47"
$ws.Range("E3").Value = "'-1.47%"
$ws.Range("D4").Value = "'5.091"
$ws.Range("D5").Value = "'0.07948"
$ws.Range("E5").Value = "'6.38%"
$ws.Range("D6").Value = "'2.193"
$ws.Range("E6").Value = "'-10.13%"
$ws.Range("D7").Value = "'7.744"
$ws.Range("E7").Value = "'-3.32%"
$ws.Range("E8").Value = "'0.17%"
$ws.Range("D9").Value = "'0.9171"
$ws.Range("E9").Value = "'-0.19%"
$ws.Range("E10").Value = "'-0.45%"
$ws.Range("D11").Value = "'0.07340"
$ws.Range("E11").Value = "'-3.52%"
$ws.Range("D12").Value = "'0.09120"
$ws.Range("D13").Value = "'0.03026"
$ws.Range("E13").Value = "'0.54%"
$ws.Range("E14").Value = "'0.80%"
$ws.Range("D15").Value = "'0.001504"
$ws.Range("E15").Value = "'-0.63%"
$ws.Range("D16").Value = "'0.005965"
$ws.Range("E16").Value = "'-1.89%"
$ws.Range("D17").Value = "'3.480"
$ws.Range("E17").Value = "'-0.76%"
$ws.Range("E18").Value = "'1.70%"
$ws.Range("E19").Value = "'0.31%"
$ws.Range("E20").Value = "'-2.14%"
$ws.Range("D21").Value = "'4.193"
$ws.Range("E21").Value = "'-9.94%"
$ws.Range("E22").Value = "'8.62%"
$ws.Range("D23").Value = "'0.04624"
$ws.Range("D24").Value = "'0.001241"
$ws.Range("E24").Value = "'-1.63%"
$ws.Range("D25").Value = "'0.004460"
$ws.Range("E25").Value = "'-1.47%"
$ws.Range("E26").Value = "'-7.74%"
$ws.Range("E27").Value = "'23.96%"
$ws.Range("D39").Value = "'0.01745"
$ws.Range("E39").Value = "'-2.21%"
$ws.Range("D40").Value = "'0.04594"
$ws.Range("E40").Value = "'0.65%"
$ws.Range("D41").Value = "'0.006938"
$ws.Range("E41").Value = "'-6.06%"
$ws.Range("D42").Value = "'0.1358"
$ws.Range("E42").Value = "'-0.30%"
$ws.Range("E43").Value = "'0.40%"
$ws.Range("D44").Value = "'0.009542"
$ws.Range("E44").Value = "'-11.59%"
$ws.Range("D45").Value = "'0.00006293"
$ws.Range("E45").Value = "'-2.42%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.06%"
$ws.Range("D47").Value = "'0.007975"
$ws.Range("E47").Value = "'-19.32%"
$ws.Range("D48").Value = "'0.7471"
$ws.Range("E48").Value = "'-8.95%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.06%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'0.01%"
